$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename to snake_case English field names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalize Spanish prepositions/articles ("de"/"el"/"la"/"del") in place names ---
$ws.Range("A15").Value  = "Ciudad De México"
$ws.Range("A23").Value  = "Estado De México"
$ws.Range("B23").Value  = "Ecatepec De Morelos"
$ws.Range("B24").Value  = "San Antonio La Isla"
$ws.Range("B25").Value  = "Tlalnepantla De Baz"
$ws.Range("B28").Value  = "Apaseo El Grande"
$ws.Range("B30").Value  = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B34").Value  = "Coyuca De Benítez"
$ws.Range("B35").Value  = "Coyuca De Catalán"
$ws.Range("B36").Value  = "Cutzamala De Pinzón"
$ws.Range("B38").Value  = "Taxco De Alarcón"
$ws.Range("B44").Value  = "Mineral Del Chico"
$ws.Range("B45").Value  = "Pachuca De Soto"
$ws.Range("B46").Value  = "Tulancingo De Bravo"
$ws.Range("B48").Value  = "Atotonilco El Alto"
$ws.Range("B50").Value  = "Encarnación De Díaz"
$ws.Range("B73").Value  = "Tlacolula De Matamoros"
$ws.Range("B82").Value  = "Cadereyta De Montes"
$ws.Range("B100").Value = "Martínez De La Torre"

# --- Drop the trailing metadata/footer rows (sample size, source, author, date) ---
$ws.Rows("112:116").Delete()
